# Add MIT license info for biogrid studies.
#
# Adds two new trailing columns to Table1 ("license_label", "license_id"),
# fills them in for the biogrid_* data rows (rows 2-11; the last data row,
# row 12, is intentionally left blank to match the source edit), and tidies
# up a couple of cosmetic details (column widths, active selection) that
# came along with the authoring session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- New table columns -----------------------------------------------------
$licenseLabelCol = $lo.ListColumns.Add()
$ws.Range("CE1").Value = "license_label"

$licenseIdCol = $lo.ListColumns.Add()
$ws.Range("CF1").Value = "license_id"

# --- Fill license values for the biogrid studies (rows 2-11) --------------
$ws.Range("CE2:CE11").Value = "MIT License"
$ws.Range("CF2:CF11").Value = "SWO:9000074"

# --- Column widths (best effort match to the authored widths) -------------
$ws.Columns.Item(82).ColumnWidth = 255
$ws.Columns.Item(83).ColumnWidth = 13.6666666666667
$ws.Columns.Item(84).ColumnWidth = 11.8333333333333

# --- Restore the active cell selection -------------------------------------
$ws.Range("CD10").Select() | Out-Null
